$d = $word.ActiveDocument

function Get-ParaByExactText($text) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

function Get-CleanTemplatePara() {
    # The empty paragraph with Times New Roman / sz24 / non-bold formatting
    # (no direct text), used as a formatting donor for new answer paragraphs.
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.TrimEnd([char]13, [char]7) -eq "" -and $p.Range.Font.Size -eq 12 -and $p.Range.Font.Bold -eq 0) {
            return $p
        }
    }
    return $null
}

function Insert-AnswerParagraphAfter($anchorPara, [string]$text) {
    # Build a fresh paragraph with the clean (Times New Roman, 12pt,
    # non-bold) formatting right after the template paragraph, fill it
    # with $text, then cut it and paste it immediately after $anchorPara.
    $tmpl = Get-CleanTemplatePara
    $tmplRange = $tmpl.Range
    $tmplRange.InsertParagraphAfter()

    $newParaIndex = $tmpl.Index + 1
    $newPara = $d.Paragraphs.Item($newParaIndex)
    $newPara.Range.InsertBefore($text)

    $cutRange = $newPara.Range
    $cutRange.Cut()

    $insertPoint = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)
    $insertPoint.Select()
    $word.Selection.Paste()
}

$p2 = Get-ParaByExactText("2-")
Insert-AnswerParagraphAfter $p2 "V"

$p2 = Get-ParaByExactText("2-")
$vPara = $p2.Next()
Insert-AnswerParagraphAfter $vPara "F"

$p2 = Get-ParaByExactText("2-")
$fPara1 = $p2.Next().Next()
Insert-AnswerParagraphAfter $fPara1 "F"

$p2 = Get-ParaByExactText("2-")
$fPara2 = $p2.Next().Next().Next()
Insert-AnswerParagraphAfter $fPara2 "Letra A"

$p3a = Get-ParaByExactText("3-a)")
Insert-AnswerParagraphAfter $p3a "Performance, automação (reutilização de código)."

$p3a = Get-ParaByExactText("3-a)")
$p3b = $p3a.Next().Next()
Insert-AnswerParagraphAfter $p3b "Manutenção, dependente da sintaxe do Banco de Dados."

Write-Output "All answer paragraphs inserted. Paragraph count: $($d.Paragraphs.Count)"
